$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The dataset is a sliding window of x/y/z sensor samples: a new reading is
# prepended at row 2, every existing reading shifts down one row, and the
# oldest reading (which fell off the bottom of the window) is dropped so the
# sheet keeps the same 20-row data window (rows 2-21).

$ws.Cells.Item(2, 1).Value = -0.0930042341351509
$ws.Cells.Item(2, 2).Value = 0.1044579595327377
$ws.Cells.Item(2, 3).Value = -0.295353353023529
$ws.Cells.Item(3, 1).Value = 0.0355829000473022
$ws.Cells.Item(3, 2).Value = 0.3640756905078888
$ws.Cells.Item(3, 3).Value = -0.0638354197144508
$ws.Cells.Item(4, 1).Value = 0.5752823352813721
$ws.Cells.Item(4, 2).Value = 0.5053382515907288
$ws.Cells.Item(4, 3).Value = -0.007177666760981
$ws.Cells.Item(5, 1).Value = 0.5288565754890442
$ws.Cells.Item(5, 2).Value = 0.7619016766548157
$ws.Cells.Item(5, 3).Value = 0.2580905556678772
$ws.Cells.Item(6, 1).Value = -0.7357872128486633
$ws.Cells.Item(6, 2).Value = -0.0797179117798805
$ws.Cells.Item(6, 3).Value = 0.8136724829673767
$ws.Cells.Item(7, 1).Value = -1.146135926246643
$ws.Cells.Item(7, 2).Value = -0.6188064813613892
$ws.Cells.Item(7, 3).Value = -0.1640173196792602
$ws.Cells.Item(8, 1).Value = 0.8213083148002625
$ws.Cells.Item(8, 2).Value = -3.046079635620117
$ws.Cells.Item(8, 3).Value = 1.09803032875061
$ws.Cells.Item(9, 1).Value = 1.518305540084839
$ws.Cells.Item(9, 2).Value = -0.5958990454673767
$ws.Cells.Item(9, 3).Value = 0.1006400510668754
$ws.Cells.Item(10, 1).Value = -0.9292787313461304
$ws.Cells.Item(10, 2).Value = 0.0415388382971286
$ws.Cells.Item(10, 3).Value = -2.335643291473389
$ws.Cells.Item(11, 1).Value = 0.5500841736793518
$ws.Cells.Item(11, 2).Value = 1.703092336654663
$ws.Cells.Item(11, 3).Value = -0.4928155243396759
$ws.Cells.Item(12, 1).Value = 0.2722931802272796
$ws.Cells.Item(12, 2).Value = -0.3020728528499603
$ws.Cells.Item(12, 3).Value = 0.2756529450416565
$ws.Cells.Item(13, 1).Value = -0.6982190012931824
$ws.Cells.Item(13, 2).Value = 0.4492913782596588
$ws.Cells.Item(13, 3).Value = -0.5900958180427551
$ws.Cells.Item(14, 1).Value = -0.4308127164840698
$ws.Cells.Item(14, 2).Value = 1.22447943687439
$ws.Cells.Item(14, 3).Value = -0.319024384021759
$ws.Cells.Item(15, 1).Value = -0.0684169083833694
$ws.Cells.Item(15, 2).Value = 0.4977024495601654
$ws.Cells.Item(15, 3).Value = -0.0363464802503585
$ws.Cells.Item(16, 1).Value = -0.1298088580369949
$ws.Cells.Item(16, 2).Value = 0.0172569435089826
$ws.Cells.Item(16, 3).Value = 0.1985312104225158
$ws.Cells.Item(17, 1).Value = 0.0598647929728031
$ws.Cells.Item(17, 2).Value = 0.5253441333770752
$ws.Cells.Item(17, 3).Value = -0.0415388382971286
$ws.Cells.Item(18, 1).Value = -0.1009454801678657
$ws.Cells.Item(18, 2).Value = 0.2724458873271942
$ws.Cells.Item(18, 3).Value = -0.4051563739776611
$ws.Cells.Item(19, 1).Value = -0.0010690141934901
$ws.Cells.Item(19, 2).Value = -0.1505782902240753
$ws.Cells.Item(19, 3).Value = 0.0145080499351024
$ws.Cells.Item(20, 1).Value = -0.1108720451593399
$ws.Cells.Item(20, 2).Value = 0.0652098655700683
$ws.Cells.Item(20, 3).Value = -0.131183311343193
$ws.Cells.Item(21, 1).Value = 0.0563523173332214
$ws.Cells.Item(21, 2).Value = 0.0316122770309448
$ws.Cells.Item(21, 3).Value = 0.1798998117446899

# Drop the row that is now beyond the window (the old last row, row 22)
$ws.Rows.Item(22).Delete()
